# Generate Report for Handback
#
# The two localization files "bed71e79-5634-4f78-be6c-f2c1dca6cf5b" and
# "0055612a-0e2e-4e43-b43c-dc47a2018e97" have effectively swapped places
# in the report: what used to be 0055612a (still "Ready for handoff") is
# now handed back in sync with en-US, with fresh handoff/handback
# timestamps, while bed71e79 keeps "Handed back" status with refreshed
# handback datetimes. This script rewrites the three worksheets (Overview,
# zh-cn, de-de) cell-by-cell to match the new report content, and then
# re-points the hyperlinks' display text at the same underlying addresses.

function Update-HyperlinkTexts($ws, $textMap) {
    # Capture existing hyperlinks (range + address) before touching cell
    # values, since editing .Value does not keep the hyperlink "display"
    # text in sync in this engine.
    $links = @()
    foreach ($hl in $ws.Hyperlinks) {
        $links += [PSCustomObject]@{ Range = $hl.Range.Address(); Address = $hl.Address }
    }

    # Remove all existing hyperlinks on the sheet.
    while ($ws.Hyperlinks.Count -gt 0) {
        foreach ($hl in $ws.Hyperlinks) {
            $hl.Delete()
            break
        }
    }

    # Re-create each hyperlink at the same range / same target address,
    # with the display text updated to match the new cell content.
    foreach ($l in $links) {
        $newText = $textMap[$l.Range]
        if (-not $newText) { $newText = $l.Range }
        $ws.Hyperlinks.Add($ws.Range($l.Range), $l.Address, "", "", $newText)
    }
}

$wb = $excel.ActiveWorkbook

# ===================== Sheet: Overview =====================
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "0055612a-0e2e-4e43-b43c-dc47a2018e97.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "2016-03-30 11:02:32"

$ws.Range("A3").Value = "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "2016-03-30 11:02:32"

$overviewMap = @{
    "`$A`$2" = "0055612a-0e2e-4e43-b43c-dc47a2018e97.md"
    "`$A`$3" = "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.md"
}
Update-HyperlinkTexts $ws $overviewMap

# ===================== Sheet: zh-cn =====================
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "0055612a-0e2e-4e43-b43c-dc47a2018e97.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "0055612a-0e2e-4e43-b43c-dc47a2018e97.3cb6669bbc1d8c1478dd6c22dae19bb2cb345e54.zh-cn.xlf"
$ws.Range("E2").Value = "2016-03-30 11:02:18"
$ws.Range("F2").Value = "0055612a-0e2e-4e43-b43c-dc47a2018e97.md"
$ws.Range("G2").Value = "0055612a-0e2e-4e43-b43c-dc47a2018e97.3cb6669bbc1d8c1478dd6c22dae19bb2cb345e54.zh-cn.xlf"
$ws.Range("H2").Value = "2016-03-30 11:03:30"
$ws.Range("J2").Value = "Include"

$ws.Range("A3").Value = "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.a3205878e3a5027496ef44cb964c4bde2ca4dc1b.zh-cn.xlf"
$ws.Range("E3").Value = "2016-03-30 11:02:18"
$ws.Range("F3").Value = "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.md"
$ws.Range("G3").Value = "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.a3205878e3a5027496ef44cb964c4bde2ca4dc1b.zh-cn.xlf"
$ws.Range("H3").Value = "2016-03-30 11:03:30"
$ws.Range("J3").Value = "Include"

$zhcnMap = @{
    "`$A`$2" = "0055612a-0e2e-4e43-b43c-dc47a2018e97.md"
    "`$D`$2" = "0055612a-0e2e-4e43-b43c-dc47a2018e97.3cb6669bbc1d8c1478dd6c22dae19bb2cb345e54.zh-cn.xlf"
    "`$F`$2" = "0055612a-0e2e-4e43-b43c-dc47a2018e97.md"
    "`$G`$2" = "0055612a-0e2e-4e43-b43c-dc47a2018e97.3cb6669bbc1d8c1478dd6c22dae19bb2cb345e54.zh-cn.xlf"
    "`$A`$3" = "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.md"
    "`$D`$3" = "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.a3205878e3a5027496ef44cb964c4bde2ca4dc1b.zh-cn.xlf"
    "`$F`$3" = "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.md"
    "`$G`$3" = "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.a3205878e3a5027496ef44cb964c4bde2ca4dc1b.zh-cn.xlf"
}
Update-HyperlinkTexts $ws $zhcnMap

# ===================== Sheet: de-de =====================
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "0055612a-0e2e-4e43-b43c-dc47a2018e97.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "0055612a-0e2e-4e43-b43c-dc47a2018e97.3cb6669bbc1d8c1478dd6c22dae19bb2cb345e54.de-de.xlf"
$ws.Range("E2").Value = "2016-03-30 11:02:32"
$ws.Range("F2").Value = "0055612a-0e2e-4e43-b43c-dc47a2018e97.md"
$ws.Range("G2").Value = "0055612a-0e2e-4e43-b43c-dc47a2018e97.3cb6669bbc1d8c1478dd6c22dae19bb2cb345e54.de-de.xlf"
$ws.Range("H2").Value = "2016-03-30 11:03:49"
$ws.Range("J2").Value = "Include"

$ws.Range("A3").Value = "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.a3205878e3a5027496ef44cb964c4bde2ca4dc1b.de-de.xlf"
$ws.Range("E3").Value = "2016-03-30 11:02:32"
$ws.Range("F3").Value = "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.md"
$ws.Range("G3").Value = "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.a3205878e3a5027496ef44cb964c4bde2ca4dc1b.de-de.xlf"
$ws.Range("H3").Value = "2016-03-30 11:03:49"
$ws.Range("J3").Value = "Include"

$dedeMap = @{
    "`$A`$2" = "0055612a-0e2e-4e43-b43c-dc47a2018e97.md"
    "`$D`$2" = "0055612a-0e2e-4e43-b43c-dc47a2018e97.3cb6669bbc1d8c1478dd6c22dae19bb2cb345e54.de-de.xlf"
    "`$F`$2" = "0055612a-0e2e-4e43-b43c-dc47a2018e97.md"
    "`$G`$2" = "0055612a-0e2e-4e43-b43c-dc47a2018e97.3cb6669bbc1d8c1478dd6c22dae19bb2cb345e54.de-de.xlf"
    "`$A`$3" = "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.md"
    "`$D`$3" = "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.a3205878e3a5027496ef44cb964c4bde2ca4dc1b.de-de.xlf"
    "`$F`$3" = "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.md"
    "`$G`$3" = "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.a3205878e3a5027496ef44cb964c4bde2ca4dc1b.de-de.xlf"
}
Update-HyperlinkTexts $ws $dedeMap

"Report regenerated for handback"
